# Rename the existing sheet to ibkGuerza, then add a new sheet "ibkIltis" after it,
# re-creating the contents of ibkGuerza's generation-log table (with a tweaked first ID).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ibkGuerza"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ibkIltis"

$ws2.Cells.Item(1,1).Value = 'Generation Log'
$ws2.Cells.Item(3,1).Value = 'idx_s'
$ws2.Cells.Item(3,2).Value = 'Description'
$ws2.Cells.Item(3,3).Value = 'Variable parameter'
$ws2.Cells.Item(3,4).Value = 'mesh factor'
$ws2.Cells.Item(3,5).Value = '# of Bridges'
$ws2.Cells.Item(3,6).Value = 'Status'
$ws2.Cells.Item(4,1).Value = 501
$ws2.Cells.Item(4,2).Value = 'Hero Parameters'
$ws2.Cells.Item(4,3).Value = 'none (single bridge)'
$ws2.Cells.Item(4,4).Value = 3
$ws2.Cells.Item(4,5).Value = 1
$ws2.Cells.Item(5,1).Value = 10
$ws2.Cells.Item(5,2).Value = 'Hero Parameters'
$ws2.Cells.Item(5,3).Value = 'none (single bridge)'
$ws2.Cells.Item(5,4).Value = 1
$ws2.Cells.Item(5,5).Value = 1
$ws2.Cells.Item(7,1).Value = 20
$ws2.Cells.Item(7,2).Value = 'Hero Parameters'
$ws2.Cells.Item(7,3).Value = 't_p (200-800)'
$ws2.Cells.Item(7,4).Value = 3
$ws2.Cells.Item(7,5).Value = 10
$ws2.Cells.Item(8,1).Value = 21
$ws2.Cells.Item(8,2).Value = 'Hero Parameters (but L=7000)'
$ws2.Cells.Item(8,3).Value = 't_p (200-800)'
$ws2.Cells.Item(8,4).Value = 3
$ws2.Cells.Item(8,5).Value = 10
$ws2.Cells.Item(9,1).Value = 22
$ws2.Cells.Item(9,2).Value = 'Hero Parameters (but L=10000)'
$ws2.Cells.Item(9,3).Value = 't_p (200-800)'
$ws2.Cells.Item(9,4).Value = 3
$ws2.Cells.Item(9,5).Value = 10
$ws2.Cells.Item(12,1).Value = 30
$ws2.Cells.Item(12,2).Value = 'Hero Parameters'
$ws2.Cells.Item(12,3).Value = 't_w (200-800)'
$ws2.Cells.Item(12,4).Value = 3
$ws2.Cells.Item(12,5).Value = 10
$ws2.Cells.Item(15,1).Value = 40
$ws2.Cells.Item(15,2).Value = 'Hero Parameters'
$ws2.Cells.Item(15,3).Value = 'L (3000,8000)'
$ws2.Cells.Item(15,4).Value = 3
$ws2.Cells.Item(15,5).Value = 10
$ws2.Cells.Item(17,1).Value = 50
$ws2.Cells.Item(17,2).Value = 'Hero Parameters'
$ws2.Cells.Item(17,3).Value = 'h_w(2000,5000)'
$ws2.Cells.Item(17,4).Value = 3
$ws2.Cells.Item(17,5).Value = 10
$ws2.Cells.Item(19,1).Value = 60
$ws2.Cells.Item(19,2).Value = 'Hero Parameters'
$ws2.Cells.Item(19,3).Value = 'b1 (10000,20000)'
$ws2.Cells.Item(19,4).Value = 3
$ws2.Cells.Item(19,5).Value = 10
$ws2.Cells.Item(21,1).Value = 70
$ws2.Cells.Item(21,2).Value = 'Hero Parameters'
$ws2.Cells.Item(21,3).Value = 'd1_plate(10,30)'
$ws2.Cells.Item(21,4).Value = 3
$ws2.Cells.Item(21,5).Value = 10
$ws2.Cells.Item(23,1).Value = 80
$ws2.Cells.Item(23,2).Value = 'Hero Parameters'
$ws2.Cells.Item(23,3).Value = 'd4_plate(10,30)'
$ws2.Cells.Item(23,4).Value = 3
$ws2.Cells.Item(23,5).Value = 10
$ws2.Cells.Item(25,1).Value = 90
$ws2.Cells.Item(25,2).Value = 'Hero Parameters'
$ws2.Cells.Item(25,3).Value = 's_plate(75,250)'
$ws2.Cells.Item(25,4).Value = 3
$ws2.Cells.Item(25,5).Value = 10
$ws2.Cells.Item(27,1).Value = 100
$ws2.Cells.Item(27,2).Value = 'Hero Parameters (L=7000)'
$ws2.Cells.Item(27,3).Value = 'd1_walls(10,30)'
$ws2.Cells.Item(27,4).Value = 3
$ws2.Cells.Item(27,5).Value = 10
$ws2.Cells.Item(29,1).Value = 110
$ws2.Cells.Item(29,2).Value = 'Hero Parameters (L=7000)'
$ws2.Cells.Item(29,3).Value = 'd4_walls(10,30)'
$ws2.Cells.Item(29,4).Value = 3
$ws2.Cells.Item(29,5).Value = 10
$ws2.Cells.Item(31,1).Value = 120
$ws2.Cells.Item(31,2).Value = 'Hero Parameters (L=7000)'
$ws2.Cells.Item(31,3).Value = 's_walls(75,250)'
$ws2.Cells.Item(31,4).Value = 3
$ws2.Cells.Item(31,5).Value = 10
$ws2.Cells.Item(33,1).Value = 130
$ws2.Cells.Item(33,2).Value = 'Hero Parameters'
$ws2.Cells.Item(33,3).Value = 'fcc (10-90)'
$ws2.Cells.Item(34,1).Value = 131
$ws2.Cells.Item(34,2).Value = 'Hero Parameters (but L=7000)'
$ws2.Cells.Item(34,3).Value = 'fcc (10-90)'
$ws2.Cells.Item(34,4).Value = 3
$ws2.Cells.Item(34,5).Value = 10
$ws2.Cells.Item(36,1).Value = 140
$ws2.Cells.Item(36,2).Value = 'Hero Parameters'
$ws2.Cells.Item(36,3).Value = 's (0,1)'
$ws2.Cells.Item(36,4).Value = 3
$ws2.Cells.Item(36,5).Value = 10
$ws2.Cells.Item(38,1).Value = 150
$ws2.Cells.Item(38,2).Value = 'Hero Parameters '
$ws2.Cells.Item(38,3).Value = 'beta (-89,89)'
$ws2.Cells.Item(38,4).Value = 3
$ws2.Cells.Item(38,5).Value = 10

# Column C on ibkIltis is widened (matches the style used on ibkGuerza's B:E columns).
$ws2.Columns.Item(3).ColumnWidth = 19.6

# Update view/selection state: ibkGuerza is scrolled back up and a data range is
# highlighted; ibkIltis becomes the active/selected tab with its default top-left view.
$ws1.Activate()
$ws1.Range("A5:E38").Select()

$ws2.Activate()
$ws2.Range("H17").Select()
